# Applies the frontendPackageOverview.pptx diagram-tidy-up edit:
#   - two "Elbow Connector" shapes are turned into plain straight connectors
#     (bentConnector3 -> straightConnector1) and re-routed to their new,
#     slightly shorter span;
#   - the three big column header textboxes ("View" / "Components" /
#     "Storage") are narrowed and re-positioned now that the boxes above them
#     moved closer together;
#   - the third textbox's caption is shortened from "Storage" to "Store".
#
# Point values below are not simply EMU/12700: Shape.Left/Top/Width/Height
# round-trips through a single-precision (float32) store, so each literal is
# chosen (mid-way through its valid float32 bucket) so it reliably lands on
# the exact target EMU value after conversion, instead of drifting +/-1 EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Elbow Connector 182": bentConnector3 -> straightConnector1 ----------
$conn1 = $s.Shapes.Item("Elbow Connector 182")
$conn1.ConnectorFormat.Type = 1   # msoConnectorStraight
$conn1.Left   = 225.91303252598425
$conn1.Top    = 272.90295408582676
$conn1.Width  = 29.50751971496063
$conn1.Height = 0

# --- "Elbow Connector 184": bentConnector3 -> straightConnector1 ----------
$conn2 = $s.Shapes.Item("Elbow Connector 184")
$conn2.ConnectorFormat.Type = 1   # msoConnectorStraight
$conn2.VerticalFlip = 0           # clear the flipV="1" left over from the bent version
$conn2.Left   = 401.74658203307087
$conn2.Top    = 275.0920868440945
$conn2.Width  = 30.27240180472441
$conn2.Height = 0

# --- "TextBox 131" ("View") header ----------------------------------------
$viewBox = $s.Shapes.Item("TextBox 131")
$viewBox.Left   = 93.58051294094489
$viewBox.Top    = 246.99743657480315
$viewBox.Width  = 108.90145492283466
$viewBox.Height = 38.77507777007874

# --- "TextBox 132" ("Components") header ----------------------------------
$componentsBox = $s.Shapes.Item("TextBox 132")
$componentsBox.Left   = 239.4968871937008
$componentsBox.Top    = 250.33067322125987
$componentsBox.Width  = 169.48980717952753
$componentsBox.Height = 36.351614003149606

# --- "TextBox 133" ("Storage" -> "Store") header --------------------------
$storeBox = $s.Shapes.Item("TextBox 133")
$storeBox.Left   = 432.4256286511811
$storeBox.Top    = 249.7664184527559
$storeBox.Width  = 146.99491878976374
$storeBox.Height = 38.77507777007874
$storeBox.TextFrame.TextRange.Text = "Store"
